$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old header row 2 (ACCESSION NO / TITLE / ACQUISITION NO. / Item / ITEM DESCRIPTION / LOCATION | SECTION / QTY / AR NUMBER / 2022/23 RFID Number / COLLECTIONS)
# this shifts rows 3-68 up to become rows 2-67
$ws.Rows.Item(2).Delete()

# Clear the acquisition-date values in column D for the first three data rows (keep the one with 1936)
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()

# Re-establish K1:N1 as visible header cells (QTY / AR NUMBER / 2022/23 RFID Number / COLLECTIONS)
$ws.Range("K1").Value = "QTY"
$ws.Range("L1").Value = "AR NUMBER"
$ws.Range("M1").Value = "2022/23 RFID Number"
$ws.Range("N1").Value = "COLLECTIONS"

# Match the bold header style used by A1:H1
$ws.Range("A1").Copy()
$ws.Range("K1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Unhide columns K:N and give them normal width, resize H/I back down from the old wide layout
$ws.Range("K1:N1").EntireColumn.Hidden = $false
$ws.Range("K1:N1").EntireColumn.ColumnWidth = 11.833333333333332
$ws.Range("H1").EntireColumn.ColumnWidth = 11.666666666666666
$ws.Range("I1").EntireColumn.ColumnWidth = 8.0

# Move the active selection back to A2
[void]$ws.Range("A2").Select()
